$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 161
    3  = 7112
    4  = 5094
    5  = 73
    9  = 96
    10 = 74
    11 = 79
    12 = 191
    13 = 623
    14 = 171
    15 = 50
    16 = 131
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
